$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "temp" (row 7) and "sal" (row 8) attribute rows, which are not
# provided in RaEn617withbiosat.csv. Deleting these rows shifts the
# following rows (biosat, O2_Ar_ratio_corrected) up by two.
$ws.Rows("7:8").Delete()

# Update the active selection to match the post-edit state (rows that were
# just removed/shifted), mirroring the selection left after the row delete.
$ws.Range("A7:XFD8").Select()
